# 6.2.1.1 — add the 2020 reporting column (M) with the new survey values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the thin separator row into the new column M (copy format from L3) ---
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

# --- Row 4: new year header 2020 in M4 (copy format from L4) ---
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 2020

# --- Row 5: Kyrgyz Republic total, 2020 value ---
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 34.377950588852634

# --- Row 6: Batken oblast, 2020 value ---
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 4.8358243107925931

# --- Row 7: Jalal-Abat oblast, 2020 value (M7 already exists, blank) ---
$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = 5.9543034993102522

# --- Row 8: Yssyk-Kul oblast, 2020 value ---
$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = 51.21106605430419

# --- Row 9: Naryn oblast, 2020 value ---
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M9").Value = 27.156801192263725

# --- Row 10: Osh oblast, 2020 value ---
$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("M10").Value = 0.94331159862228353

# --- Row 11: Talas oblast, 2020 value ---
$ws.Range("L11").Copy()
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("M11").Value = 7.8509592890793316

# --- Row 12: Chui oblast, 2020 value ---
$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("M12").Value = 64.733302669743793

# --- Row 13: Bishkek city, 2020 value ---
$ws.Range("L13").Copy()
$ws.Range("M13").PasteSpecial(-4122)
$ws.Range("M13").Value = 97.67954817102779

# --- Row 14: Osh city, 2020 value ---
$ws.Range("L14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = 46.725153243037099

# --- Move the active selection to L19 (matches the saved cursor position) ---
[void]$ws.Range("L19").Select()
